# Daily refresh of the cryptos price/volume list.
# Price ("D") and Volume(1h) ("E") columns are updated in place for each
# coin row; rows 27/28 additionally swap their Coin/Link/Price/Volume
# content (InternetComputer(DFINITY) and Dai traded ranks), while the
# row index in column A stays put.
#
# Column D stores plain text look-alikes of numbers (e.g. "682.11",
# "69.399.99"), so NumberFormat is forced to Text ("@") before writing
# each value - otherwise Excel's COM layer would auto-coerce them to
# numeric cells (dropping formatting such as trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.399.99'
$ws.Range('E2').Value = '  -1.82%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.686.12'
$ws.Range('E3').Value = '  -2.90%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '682.11'
$ws.Range('E5').Value = '  -2.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.49'
$ws.Range('E6').Value = '  -3.95%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.684.75'
$ws.Range('E7').Value = '  -2.89%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -4.05%  '
$ws.Range('E10').Value = '  -7.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.34'
$ws.Range('E11').Value = '  -3.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.446'
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('E13').Value = '  -4.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.46'
$ws.Range('E14').Value = '  -6.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.309.17'
$ws.Range('E15').Value = '  -2.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.687.54'
$ws.Range('E16').Value = '  -3.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.371.50'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.25'
$ws.Range('E19').Value = '  -6.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.62'
$ws.Range('E20').Value = '  -6.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '482.51'
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.87'
$ws.Range('E22').Value = '  -7.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.666'
$ws.Range('E23').Value = '  -7.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.28'
$ws.Range('E24').Value = '  -4.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.832.99'
$ws.Range('E25').Value = '  -3.00%  '
$ws.Range('E26').Value = '  -8.34%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.49'
$ws.Range('E28').Value = '  -4.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.61'
$ws.Range('E29').Value = '  -6.13%  '
$ws.Range('E30').Value = '  -8.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.73'
$ws.Range('E31').Value = '  -9.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.08'
$ws.Range('E32').Value = '  -8.11%  '
$ws.Range('E33').Value = '  -6.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.07'
$ws.Range('E34').Value = '  -6.42%  '
$ws.Range('E35').Value = '  -4.52%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.654.24'
$ws.Range('E37').Value = '  -3.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.50'
$ws.Range('E38').Value = '  -5.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.35'
$ws.Range('E39').Value = '  +7.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0934'
$ws.Range('E40').Value = '  -7.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.25'
$ws.Range('E41').Value = '  -4.18%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('E44').Value = '  -6.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '160.47'
$ws.Range('E45').Value = '  -3.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '48.40'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.87'
$ws.Range('E47').Value = '  -10.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.36'
$ws.Range('E48').Value = '  +1.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000289'
$ws.Range('E49').Value = '  -8.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '29.72'
$ws.Range('E50').Value = '  +5.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '396.35'
$ws.Range('E51').Value = '  -5.89%  '
